$d = $word.ActiveDocument

# The document currently ends with an empty paragraph (after the Athena
# paragraph). Append three new paragraphs after it: "New commit:", the
# Eros paragraph, and a trailing empty paragraph.

$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$lastPara.Range.InsertParagraphAfter()

$commitPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$commitPara.Range.InsertAfter('New commit:')

$commitPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$commitPara.Range.InsertParagraphAfter()

$erosPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$erosPara.Range.InsertAfter('n Greek mythology, Eros (UK: /ˈɪərɒs, ˈɛrɒs/, US: /ˈɛrɒs, ˈɛroʊs/;[2] Ancient Greek: Ἔρως, romanized: Érōs, lit. ''Desire'') is the Greek god of love and sex. His Roman counterpart was Cupid ("desire").[3] In the earliest account, he is a primordial god, while in later accounts he is described as one of the children of Aphrodite and Ares and, with some of his siblings, was one of the Erotes, a group of winged love gods')

$erosPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$erosPara.Range.InsertParagraphAfter()
